$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 667.375
$ws.Range("I2").Value = 548
$ws.Range("J2").Value = 739
$ws.Range("K2").Value = 548
$ws.Range("L2").Value = 739
$ws.Range("M2").Value = -435
$ws.Range("N2").Value = -965

$ws.Range("H12").Value = 1176.7778
$ws.Range("I12").Value = 1176.7778
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 1176.7778
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -1006.7778
$ws.Range("N12").ClearContents()

$ws.Range("H40").Value = 2520
$ws.Range("I40").Value = 2450
$ws.Range("K40").Value = 2450
$ws.Range("M40").Value = -2275

$ws.Range("H53").Value = 182.85715
$ws.Range("J53").Value = 287.75
$ws.Range("L53").Value = 287.75
$ws.Range("N53").Value = -1561.75

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 5000
$ws.Range("K62").Value = 5000
$ws.Range("M62").Value = -4376

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 5000
$ws.Range("K65").Value = 25000
$ws.Range("M65").Value = -21880

$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()

$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()

$ws.Range("H94").Value = 11332.667
$ws.Range("I94").Value = 11332.667
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 11332.667
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -10881.667
$ws.Range("N94").ClearContents()

$ws.Range("H132").Value = 615.7273
$ws.Range("I132").Value = 628.381
$ws.Range("J132").Value = 350
$ws.Range("K132").Value = 1885.143
$ws.Range("L132").Value = 1050
$ws.Range("M132").Value = 644.857
$ws.Range("N132").Value = -6110

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5101.6875
$ws.Range("I32").Value = 5101.6875
$ws.Range("K32").Value = 5101.6875
$ws.Range("M32").Value = -4814.6875

$ws.Range("H45").Value = 2285.6667
$ws.Range("I45").Value = 2285.6667
$ws.Range("K45").Value = 2285.6667
$ws.Range("M45").Value = -1908.6667

$ws.Range("H74").Value = 1332.6666
$ws.Range("I74").Value = 1332.6666
$ws.Range("K74").Value = 1332.6666
$ws.Range("M74").Value = -458.6666

$ws.Range("H77").Value = 1332.6666
$ws.Range("I77").Value = 1332.6666
$ws.Range("K77").Value = 6663.333000000001
$ws.Range("M77").Value = -2295.333000000001

$ws.Range("H96").Value = 40344
$ws.Range("J96").Value = 40344
$ws.Range("L96").Value = 40344
$ws.Range("N96").Value = -45836

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 899.125
$ws.Range("I22").Value = 884.8570999999999
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 884.8570999999999
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -711.8570999999999
$ws.Range("N22").Value = -1345

$ws.Range("H36").Value = 12916.667
$ws.Range("I36").Value = 12916.667
$ws.Range("K36").Value = 12916.667
$ws.Range("M36").Value = -12382.667

$ws.Range("H61").Value = 12250
$ws.Range("J61").Value = 12250
$ws.Range("L61").Value = 12250
$ws.Range("N61").Value = -12876

$ws.Range("H64").Value = 2000
$ws.Range("I64").Value = 2000
$ws.Range("K64").Value = 2000
$ws.Range("M64").Value = -1775

$ws.Range("H67").Value = 2000
$ws.Range("I67").Value = 2000
$ws.Range("K67").Value = 2000
$ws.Range("M67").Value = -1220

$ws.Range("H105").Value = 608.5
$ws.Range("I105").Value = 623.8889
$ws.Range("K105").Value = 623.8889
$ws.Range("M105").Value = 1123.1111

$ws.Range("H134").Value = 3425
$ws.Range("I134").Value = 3233.3333
$ws.Range("K134").Value = 9699.999899999999
$ws.Range("M134").Value = -7164.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 2749.75
$ws.Range("I10").Value = 1999.5
$ws.Range("K10").Value = 1999.5
$ws.Range("M10").Value = -1860.5

$ws.Range("H31").Value = 3826.5715
$ws.Range("I31").Value = 2881
$ws.Range("J31").Value = 9500
$ws.Range("K31").Value = 2881
$ws.Range("L31").Value = 9500
$ws.Range("M31").Value = -2586
$ws.Range("N31").Value = -10090

$ws.Range("H34").Value = 3826.5715
$ws.Range("I34").Value = 2881
$ws.Range("J34").Value = 9500
$ws.Range("K34").Value = 2881
$ws.Range("L34").Value = 9500
$ws.Range("M34").Value = -2679
$ws.Range("N34").Value = -9904

$ws.Range("H51").Value = 7666.6665
$ws.Range("I51").Value = 7666.6665
$ws.Range("K51").Value = 7666.6665
$ws.Range("M51").Value = -6930.6665

$ws.Range("H61").Value = 7666.6665
$ws.Range("I61").Value = 7666.6665
$ws.Range("K61").Value = 7666.6665
$ws.Range("M61").Value = -7318.6665

$ws.Range("H122").Value = 3075.75
$ws.Range("I122").Value = 978.4545000000001
$ws.Range("K122").Value = 2935.3635
$ws.Range("M122").Value = -485.3635000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 165.5
$ws.Range("J25").Value = 225
$ws.Range("L25").Value = 675
$ws.Range("N25").Value = -1013

$ws.Range("H26").Value = 747.5
$ws.Range("I26").Value = 346.66666
$ws.Range("K26").Value = 1039.99998
$ws.Range("M26").Value = -751.9999800000001

$ws.Range("H30").Value = 165.5
$ws.Range("J30").Value = 225
$ws.Range("L30").Value = 675
$ws.Range("N30").Value = -879

$ws.Range("H107").Value = 557.8461
$ws.Range("I107").Value = 361.6
$ws.Range("J107").Value = 680.5
$ws.Range("K107").Value = 1084.8
$ws.Range("L107").Value = 2041.5
$ws.Range("M107").Value = 835.1999999999998
$ws.Range("N107").Value = -5881.5

$ws.Range("H133").Value = 4999.5
$ws.Range("I133").Value = 4999.5
$ws.Range("K133").Value = 14998.5
$ws.Range("M133").Value = -9938.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 50000
$ws.Range("J26").Value = 50000
$ws.Range("L26").Value = 50000
$ws.Range("N26").Value = -50560

$ws.Range("H50").Value = 50000
$ws.Range("J50").Value = 50000
$ws.Range("L50").Value = 50000
$ws.Range("N50").Value = -50996

$ws.Range("H122").Value = 6687.846
$ws.Range("I122").Value = 5408.8335
$ws.Range("J122").Value = 7784.143
$ws.Range("K122").Value = 16226.5005
$ws.Range("L122").Value = 23352.429
$ws.Range("M122").Value = -13776.5005
$ws.Range("N122").Value = -28252.429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1812.3334
$ws.Range("I22").Value = 1968.75
$ws.Range("K22").Value = 1968.75
$ws.Range("M22").Value = -1673.75

$ws.Range("H27").Value = 1812.3334
$ws.Range("I27").Value = 1968.75
$ws.Range("K27").Value = 1968.75
$ws.Range("M27").Value = -1861.75

$ws.Range("H46").Value = 498.75
$ws.Range("I46").Value = 447.5
$ws.Range("J46").Value = 550
$ws.Range("K46").Value = 447.5
$ws.Range("L46").Value = 550
$ws.Range("M46").Value = -259.5
$ws.Range("N46").Value = -926

$ws.Range("H55").Value = 510.5625
$ws.Range("J55").Value = 866
$ws.Range("L55").Value = 866
$ws.Range("N55").Value = -1212

$ws.Range("H63").Value = 30000
$ws.Range("I63").Value = 30000
$ws.Range("K63").Value = 30000
$ws.Range("M63").Value = -29251

$ws.Range("H66").Value = 30000
$ws.Range("I66").Value = 30000
$ws.Range("K66").Value = 90000
$ws.Range("M66").Value = -86256

$ws.Range("H74").Value = 15000
$ws.Range("I74").Value = 15000
$ws.Range("K74").Value = 15000
$ws.Range("M74").Value = -14002

$ws.Range("H77").Value = 15000
$ws.Range("I77").Value = 15000
$ws.Range("K77").Value = 45000
$ws.Range("M77").Value = -40008

$ws.Range("H122").Value = 4100
$ws.Range("I122").Value = 3600
$ws.Range("K122").Value = 10800
$ws.Range("M122").Value = -8350

$ws.Range("H132").Value = 2816.75
$ws.Range("I132").Value = 2252
$ws.Range("K132").Value = 6756
$ws.Range("M132").Value = -4226

$ws.Range("H136").Value = 4499.75
$ws.Range("I136").Value = 3999.6667
$ws.Range("K136").Value = 11999.0001
$ws.Range("M136").Value = -9449.000100000001

$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("M140").ClearContents()
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 952.5
$ws.Range("I9").Value = 952.5
$ws.Range("K9").Value = 952.5
$ws.Range("M9").Value = -812.5

$ws.Range("H136").Value = 2327.6
$ws.Range("I136").Value = 2534.5
$ws.Range("K136").Value = 7603.5
$ws.Range("M136").Value = -5053.5
